$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update graded criteria values (column G) per the grading pass
$ws.Range("G4").Value = 4

$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("G12").Value = 8

$ws.Range("G16").Value = 0
$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 1

# These four criteria were not graded (cleared)
$ws.Range("G24").ClearContents()
$ws.Range("G25").ClearContents()
$ws.Range("G26").ClearContents()
$ws.Range("G27").ClearContents()

$ws.Range("G28").Value = 3

# Move the active selection to reflect where grading finished
$ws.Range("G29").Select()
